$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray date that used to live in A101 (data for row 101 moved down)
$ws.Range("A101").Clear()

# Row 103 - Security / Annas payment
$ws.Range("A103").Value2 = 45385
$ws.Range("B103").Value = "Security"
$ws.Range("C103").Value = "Annas"
$ws.Range("D103").Value2 = 45292
$ws.Range("E103").Value2 = 45384
$ws.Range("F103").Value = 789000
$ws.Range("G103").Value = 789000
$ws.Range("H103").Value = 93000
$ws.Range("J103").Value = 300000
$ws.Range("K103").Formula = "=G103+H103-J103"

# Row 104 - Manggi payment
$ws.Range("A104").Value2 = 45386
$ws.Range("B104").Value = "Manggi"
$ws.Range("D104").Value2 = 45352
$ws.Range("E104").Value2 = 45384
$ws.Range("F104").Value = 27678000
$ws.Range("G104").Value = 27678000
$ws.Range("J104").Value = 20000000
$ws.Range("K104").Formula = "=G104+H104-J104"

# Row 105 - Subadi payment
$ws.Range("B105").Value = "Subadi"
$ws.Range("D105").Value2 = 45360
$ws.Range("E105").Value2 = 45384
$ws.Range("F105").Value = 35306000
$ws.Range("G105").Value = 35306000
$ws.Range("H105").Formula = "=K53+K83"
$ws.Range("J105").Value = 30230000
$ws.Range("K105").Formula = "=G105+H105-J105"

# Update the view state to match where the user left off editing
$excel.ActiveWindow.ScrollRow = 94
$ws.Range("B99").Select()
